$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")

# --- Update Recommandations sheet data rows 2-28 ---
$ws1.Cells.Item(2, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 4
$ws1.Cells.Item(2, 4).Value = 719.88
$ws1.Cells.Item(2, 5).Value = 180.49
$ws1.Cells.Item(2, 6).Value = "🟡 Observer"
$ws1.Cells.Item(2, 7).Value = "➖ Neutre"

$ws1.Cells.Item(3, 1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(3, 2).Value = 0
$ws1.Cells.Item(3, 3).Value = 4
$ws1.Cells.Item(3, 4).Value = 610.62
$ws1.Cells.Item(3, 5).Value = 154.48
$ws1.Cells.Item(3, 6).Value = "🟡 Observer"
$ws1.Cells.Item(3, 7).Value = "➖ Neutre"

$ws1.Cells.Item(4, 1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(4, 2).Value = 0
$ws1.Cells.Item(4, 3).Value = 4
$ws1.Cells.Item(4, 4).Value = 596.88
$ws1.Cells.Item(4, 5).Value = 153.59
$ws1.Cells.Item(4, 6).Value = "🟡 Observer"
$ws1.Cells.Item(4, 7).Value = "➖ Neutre"

$ws1.Cells.Item(5, 1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(5, 2).Value = 0
$ws1.Cells.Item(5, 3).Value = 4
$ws1.Cells.Item(5, 4).Value = 584.01
$ws1.Cells.Item(5, 5).Value = 146.93
$ws1.Cells.Item(5, 6).Value = "🟡 Observer"
$ws1.Cells.Item(5, 7).Value = "➖ Neutre"

$ws1.Cells.Item(6, 1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(6, 2).Value = 0
$ws1.Cells.Item(6, 3).Value = 4
$ws1.Cells.Item(6, 4).Value = 471.26
$ws1.Cells.Item(6, 5).Value = 118.33
$ws1.Cells.Item(6, 6).Value = "🟡 Observer"
$ws1.Cells.Item(6, 7).Value = "➖ Neutre"

$ws1.Cells.Item(7, 1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(7, 2).Value = 0
$ws1.Cells.Item(7, 3).Value = 4
$ws1.Cells.Item(7, 4).Value = 464.2
$ws1.Cells.Item(7, 5).Value = 116.44
$ws1.Cells.Item(7, 6).Value = "🟡 Observer"
$ws1.Cells.Item(7, 7).Value = "➖ Neutre"

$ws1.Cells.Item(8, 1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(8, 2).Value = 0
$ws1.Cells.Item(8, 3).Value = 4
$ws1.Cells.Item(8, 4).Value = 376.99
$ws1.Cells.Item(8, 5).Value = 94.31
$ws1.Cells.Item(8, 6).Value = "🟡 Observer"
$ws1.Cells.Item(8, 7).Value = "➖ Neutre"

$ws1.Cells.Item(9, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Cells.Item(9, 2).Value = 3
$ws1.Cells.Item(9, 3).Value = 1
$ws1.Cells.Item(9, 4).Value = 19.19
$ws1.Cells.Item(9, 5).Value = -1.98
$ws1.Cells.Item(9, 6).Value = "🟢 Achat"
$ws1.Cells.Item(9, 7).Value = "✅ Renforcer"

$ws1.Cells.Item(10, 1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Cells.Item(10, 2).Value = 2
$ws1.Cells.Item(10, 3).Value = 0
$ws1.Cells.Item(10, 4).Value = 10.36
$ws1.Cells.Item(10, 5).Value = 7.5
$ws1.Cells.Item(10, 6).Value = "🟡 Observer"
$ws1.Cells.Item(10, 7).Value = "➖ Neutre"

$ws1.Cells.Item(11, 1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(11, 2).Value = 2
$ws1.Cells.Item(11, 3).Value = 1
$ws1.Cells.Item(11, 4).Value = 9.01
$ws1.Cells.Item(11, 5).Value = -3.79
$ws1.Cells.Item(11, 6).Value = "🟡 Observer"
$ws1.Cells.Item(11, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(12, 1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(12, 2).Value = 1
$ws1.Cells.Item(12, 3).Value = 0
$ws1.Cells.Item(12, 4).Value = 7.19
$ws1.Cells.Item(12, 5).Value = 7.19
$ws1.Cells.Item(12, 6).Value = "🟡 Observer"
$ws1.Cells.Item(12, 7).Value = "➖ Neutre"

$ws1.Cells.Item(13, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(13, 2).Value = 2
$ws1.Cells.Item(13, 3).Value = 0
$ws1.Cells.Item(13, 4).Value = 6.98
$ws1.Cells.Item(13, 5).Value = 4.4
$ws1.Cells.Item(13, 6).Value = "🟡 Observer"
$ws1.Cells.Item(13, 7).Value = "➖ Neutre"

$ws1.Cells.Item(14, 1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(14, 2).Value = 1
$ws1.Cells.Item(14, 3).Value = 0
$ws1.Cells.Item(14, 4).Value = 5
$ws1.Cells.Item(14, 5).Value = 5
$ws1.Cells.Item(14, 6).Value = "🟡 Observer"
$ws1.Cells.Item(14, 7).Value = "➖ Neutre"

$ws1.Cells.Item(15, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(15, 2).Value = 2
$ws1.Cells.Item(15, 3).Value = 1
$ws1.Cells.Item(15, 4).Value = 4.53
$ws1.Cells.Item(15, 5).Value = 4.35
$ws1.Cells.Item(15, 6).Value = "🟡 Observer"
$ws1.Cells.Item(15, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(16, 1).Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Cells.Item(16, 2).Value = 1
$ws1.Cells.Item(16, 3).Value = 0
$ws1.Cells.Item(16, 4).Value = 3.74
$ws1.Cells.Item(16, 5).Value = 3.74
$ws1.Cells.Item(16, 6).Value = "🟡 Observer"
$ws1.Cells.Item(16, 7).Value = "➖ Neutre"

$ws1.Cells.Item(17, 1).Value = "SMB CI (SMBC)"
$ws1.Cells.Item(17, 2).Value = 1
$ws1.Cells.Item(17, 3).Value = 1
$ws1.Cells.Item(17, 4).Value = 1.3
$ws1.Cells.Item(17, 5).Value = 3.45
$ws1.Cells.Item(17, 6).Value = "🟡 Observer"
$ws1.Cells.Item(17, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(18, 1).Value = "SUCRIVOIRE (SCRC)"
$ws1.Cells.Item(18, 2).Value = 1
$ws1.Cells.Item(18, 3).Value = 1
$ws1.Cells.Item(18, 4).Value = 0.88
$ws1.Cells.Item(18, 5).Value = -6.32
$ws1.Cells.Item(18, 6).Value = "🟡 Observer"
$ws1.Cells.Item(18, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(19, 1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Cells.Item(19, 2).Value = 1
$ws1.Cells.Item(19, 3).Value = 1
$ws1.Cells.Item(19, 4).Value = 0.01
$ws1.Cells.Item(19, 5).Value = 3.13
$ws1.Cells.Item(19, 6).Value = "🟡 Observer"
$ws1.Cells.Item(19, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(20, 1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(20, 2).Value = 1
$ws1.Cells.Item(20, 3).Value = 2
$ws1.Cells.Item(20, 4).Value = -0.85
$ws1.Cells.Item(20, 5).Value = -1.67
$ws1.Cells.Item(20, 6).Value = "🟡 Observer"
$ws1.Cells.Item(20, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(21, 1).Value = "UNILEVER CI (UNLC)"
$ws1.Cells.Item(21, 2).Value = 1
$ws1.Cells.Item(21, 3).Value = 1
$ws1.Cells.Item(21, 4).Value = -1.12
$ws1.Cells.Item(21, 5).Value = 6.36
$ws1.Cells.Item(21, 6).Value = "🟡 Observer"
$ws1.Cells.Item(21, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(22, 1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(22, 2).Value = 0
$ws1.Cells.Item(22, 3).Value = 1
$ws1.Cells.Item(22, 4).Value = -1.31
$ws1.Cells.Item(22, 5).Value = -1.31
$ws1.Cells.Item(22, 6).Value = "🟡 Observer"
$ws1.Cells.Item(22, 7).Value = "➖ Neutre"

$ws1.Cells.Item(23, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Cells.Item(23, 2).Value = 1
$ws1.Cells.Item(23, 3).Value = 2
$ws1.Cells.Item(23, 4).Value = -1.32
$ws1.Cells.Item(23, 5).Value = -1.38
$ws1.Cells.Item(23, 6).Value = "🟡 Observer"
$ws1.Cells.Item(23, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(24, 1).Value = "CIE CI (CIEC)"
$ws1.Cells.Item(24, 2).Value = 0
$ws1.Cells.Item(24, 3).Value = 1
$ws1.Cells.Item(24, 4).Value = -1.85
$ws1.Cells.Item(24, 5).Value = -1.85
$ws1.Cells.Item(24, 6).Value = "🟡 Observer"
$ws1.Cells.Item(24, 7).Value = "➖ Neutre"

$ws1.Cells.Item(25, 1).Value = "NESTLE CI (NTLC)"
$ws1.Cells.Item(25, 2).Value = 0
$ws1.Cells.Item(25, 3).Value = 1
$ws1.Cells.Item(25, 4).Value = -2.59
$ws1.Cells.Item(25, 5).Value = -2.59
$ws1.Cells.Item(25, 6).Value = "🟡 Observer"
$ws1.Cells.Item(25, 7).Value = "➖ Neutre"

$ws1.Cells.Item(26, 1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(26, 2).Value = 0
$ws1.Cells.Item(26, 3).Value = 2
$ws1.Cells.Item(26, 4).Value = -4.02
$ws1.Cells.Item(26, 5).Value = -1.58
$ws1.Cells.Item(26, 6).Value = "🟡 Observer"
$ws1.Cells.Item(26, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(27, 1).Value = "SICOR CI (SICC)"
$ws1.Cells.Item(27, 2).Value = 0
$ws1.Cells.Item(27, 3).Value = 1
$ws1.Cells.Item(27, 4).Value = -5.26
$ws1.Cells.Item(27, 5).Value = -5.26
$ws1.Cells.Item(27, 6).Value = "🟡 Observer"
$ws1.Cells.Item(27, 7).Value = "➖ Neutre"

$ws1.Cells.Item(28, 1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(28, 2).Value = 0
$ws1.Cells.Item(28, 3).Value = 3
$ws1.Cells.Item(28, 4).Value = -15.43
$ws1.Cells.Item(28, 5).Value = -2.6
$ws1.Cells.Item(28, 6).Value = "🔴 Vente"
$ws1.Cells.Item(28, 7).Value = "⚠️ Risque de décrochage"

# --- Remove now-obsolete trailing rows 29-31 ---
$ws1.Range("A29:G31").Delete()

# --- Update Top_YTD sheet ---
$ws2 = $wb.Worksheets.Item("Top_YTD")
$ws2.Cells.Item(2, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Cells.Item(2, 2).Value = 6043.48
$ws2.Cells.Item(3, 1).Value = "BRVM - SERVICES FINANCIERS"
$ws2.Cells.Item(3, 2).Value = 3974.68
$ws2.Cells.Item(4, 1).Value = "BRVM - INDUSTRIELS"
$ws2.Cells.Item(4, 2).Value = 3756.73
$ws2.Cells.Item(5, 1).Value = "BRVM-PRESTIGE"
$ws2.Cells.Item(5, 2).Value = 3562.29
$ws2.Cells.Item(6, 1).Value = "BRVM - SERVICES PUBLICS"
$ws2.Cells.Item(6, 2).Value = 2150.84
$ws2.Cells.Item(7, 1).Value = "BRVM - ENERGIE"
$ws2.Cells.Item(7, 2).Value = 2078.78
$ws2.Cells.Item(8, 1).Value = "BRVM - TELECOMMUNICATIONS"
$ws2.Cells.Item(8, 2).Value = 1323.71
$ws2.Cells.Item(9, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws2.Cells.Item(9, 2).Value = 20.27
$ws2.Cells.Item(10, 1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws2.Cells.Item(10, 2).Value = 10.57
$ws2.Cells.Item(11, 1).Value = "UNIWAX CI (UNXC)"
$ws2.Cells.Item(11, 2).Value = 8.91
